$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap the order of "Canada" / "Ucrania" in row 27 / 28 (country name column A) ---
# Row 27 becomes Ucrania, row 28 becomes Canada (their numeric stats also get refreshed below).
$ws.Cells.Item(27, 1).Value = "Ucrania"
$ws.Cells.Item(28, 1).Value = "Canada"

# --- Update the "last updated" timestamp string (title cell, row 1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 5 de Septiembre de 2020 a las 09:05"

function Set-CountryRow($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- Refresh the per-country statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) for the updated rows ---

# Row 27: Ucrania (new data)
Set-CountryRow 27 133787 2836 61649 69327 0 50 2811

# Row 28: Canada (data that used to belong to the old Canada row, now shifted down)
Set-CountryRow 28 131124 0 115926 6057 0 0 9141

# Row 60: Armenia
Set-CountryRow 60 44649 188 39823 3931 0 4 895

# Row 73: Australia
Set-CountryRow 73 26208 72 22331 3129 0 11 748

# Row 104: Hungria
Set-CountryRow 104 7892 510 3952 3316 0 3 624

# Row 152: Georgia
Set-CountryRow 152 1621 25 1302 300 0 0 19

# Row 175: Taiwan
Set-CountryRow 175 492 2 473 12 0 0 7
